$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Data edits -------------------------------------------------------
# The menu used to encode a "selected" and a "not selected" (free) option
# together in one comma-joined string (e.g. "Ice,No Ice" / "0.50,0.00").
# Blocking/availability is now handled elsewhere, so the sheets only need
# to keep the single real, chargeable option.

# Sheet1 "Ah Beng Drink": Option_1 / Amount_1 columns (E/F)
$ws1.Range("E2").Value = "Ice"
$ws1.Range("F2").Value = 0.5
$ws1.Range("E3").Value = "Ice"
$ws1.Range("F3").Value = 0.5
$ws1.Range("E4").Value = "Ice"
$ws1.Range("F4").Value = 0.5

# Sheet2 "Ah Lian Food": Option_1 / Amount_1 columns (E/F)
$ws2.Range("E2").Value = "Chilli"
$ws2.Range("F2").Value = "0.50"
$ws2.Range("E3").Value = "Chilli"
$ws2.Range("F3").Value = "0.50"
$ws2.Range("E4").Value = "Chilli"
$ws2.Range("F4").Value = "0.50"
$ws2.Range("E5").Value = "Ice"
$ws2.Range("F5").Value = "0.50"

# The lone left-aligned style (used by E5) now aligns left instead of right.
$ws2.Range("E5").HorizontalAlignment = -4131

# --- Column sizing on sheet1 (content got shorter/numeric) -------------
$ws1.Columns.Item(5).EntireColumn.AutoFit()
$ws1.Columns.Item(6).EntireColumn.AutoFit()

# --- Page setup on sheet1 ------------------------------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Selection / active-sheet bookkeeping -------------------------------
# Final state: "Ah Lian Food" (sheet2) is the active/selected tab with
# K11 selected; "Ah Beng Drink" (sheet1) keeps F2:F4 selected in the
# background.
$ws1.Range("F2:F4").Select()
$ws2.Activate()
$ws2.Range("K11").Select()
